$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.326.30'
$ws.Range("E2").Value = '  -0.27%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.791.81'
$ws.Range("E3").Value = '  -0.80%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.10'
$ws.Range("E5").Value = '  -0.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.18%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5330'
$ws.Range("E7").Value = '  -1.65%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3762'
$ws.Range("E8").Value = '  -1.64%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07494'
$ws.Range("E9").Value = '  -0.91%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.84'
$ws.Range("E10").Value = '  -2.80%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.095'
$ws.Range("E11").Value = '  -2.29%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("E12").Value = '  +0.02%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.50'
$ws.Range("E13").Value = '  -3.30%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.106'
$ws.Range("E14").Value = '  -1.34%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.808.96'
$ws.Range("E15").Value = '  -0.03%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.263'
$ws.Range("E16").Value = '  -1.15%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '89.09'
$ws.Range("E17").Value = '  -2.46%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001058'
$ws.Range("E18").Value = '  -1.13%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06490'
$ws.Range("E19").Value = '  +0.75%  '
$ws.Range("E20").Value = '  +0.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.27'
$ws.Range("E21").Value = '  +0.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.966'
$ws.Range("E22").Value = '  -0.21%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.361.34'
$ws.Range("E23").Value = '  -0.27%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.03'
$ws.Range("E24").Value = '  -1.45%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.082'
$ws.Range("E25").Value = '  -3.54%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.82'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.33'
$ws.Range("E27").Value = '  -1.61%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.001.80'
$ws.Range("E28").Value = '  -0.76%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.274'
$ws.Range("E29").Value = '  -6.79%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '122.11'
$ws.Range("E30").Value = '  -1.67%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.097'
$ws.Range("E31").Value = '  -4.16%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1044'
$ws.Range("E32").Value = '  +3.04%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.659'
$ws.Range("E33").Value = '  -0.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.566'
$ws.Range("E34").Value = '  -2.58%  '
$ws.Range("B35").Value = 'Algorand'
$ws.Range("C35").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.2258'
$ws.Range("E35").Value = '  -1.17%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06514'
$ws.Range("E36").Value = '  +4.48%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02279'
$ws.Range("E37").Value = '  -1.73%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.992'
$ws.Range("E38").Value = '  -0.38%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.489'
$ws.Range("E39").Value = '  -4.36%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.452'
$ws.Range("E40").Value = '  +5.17%  '
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6137'
$ws.Range("E41").Value = '  -3.30%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.189'
$ws.Range("E42").Value = '  +3.18%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.06'
$ws.Range("E43").Value = '  -4.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  +0.21%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.35'
$ws.Range("E45").Value = '  -0.23%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.678'
$ws.Range("E46").Value = '  -0.22%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5755'
$ws.Range("E47").Value = '  -3.33%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.62'
$ws.Range("E48").Value = '  +1.72%  '
$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.190'
$ws.Range("E49").Value = '  +3.73%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.924'
$ws.Range("E50").Value = '  -2.34%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06852'
$ws.Range("E51").Value = '  -0.85%  '
